# Update the two-digit multiplication practice table with newly generated problems.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="64×81=5184"},
    @{Row=1;  Col=2; Text="19×14=266"},
    @{Row=1;  Col=3; Text="59×80=4720"},
    @{Row=1;  Col=4; Text="65×25=1625"},
    @{Row=1;  Col=5; Text="17×96=1632"},

    @{Row=5;  Col=1; Text="15×36=540"},
    @{Row=5;  Col=2; Text="61×49=2989"},
    @{Row=5;  Col=3; Text="52×82=4264"},
    @{Row=5;  Col=4; Text="42×30=1260"},
    @{Row=5;  Col=5; Text="20×11=220"},

    @{Row=10; Col=1; Text="62×46=2852"},
    @{Row=10; Col=2; Text="13×59=767"},
    @{Row=10; Col=3; Text="84×92=7728"},
    @{Row=10; Col=4; Text="98×45=4410"},
    @{Row=10; Col=5; Text="57×39=2223"},

    @{Row=15; Col=1; Text="59×57=3363"},
    @{Row=15; Col=2; Text="80×65=5200"},
    @{Row=15; Col=3; Text="20×93=1860"},
    @{Row=15; Col=4; Text="22×61=1342"},
    @{Row=15; Col=5; Text="67×26=1742"},

    @{Row=20; Col=1; Text="14×21=294"},
    @{Row=20; Col=2; Text="80×97=7760"},
    @{Row=20; Col=3; Text="91×51=4641"},
    @{Row=20; Col=4; Text="94×84=7896"},
    @{Row=20; Col=5; Text="28×46=1288"}
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.Text
}
